$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.492.82"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.839.30"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.22"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5364"
$ws.Range("E7").Value = "  +2.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2911"
$ws.Range("E8").Value = "  -9.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07058"
$ws.Range("E9").Value = "  +3.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.22"
$ws.Range("E10").Value = "  -8.79%  "
$ws.Range("D11").Value = "1.840.56"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7197"
$ws.Range("E12").Value = "  -7.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07135"
$ws.Range("E13").Value = "  -8.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.05"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.979"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.66"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007914"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "26.504.62"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "2.070.96"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.585"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.980"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.127"
$ws.Range("E24").Value = "  -3.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.81"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.707"
$ws.Range("E26").Value = "  +1.91%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.097"
$ws.Range("E27").Value = "  -3.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.96"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.09"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.248"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08867"
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.017"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04788"
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7205"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.129"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.093"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.253"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01701"
$ws.Range("E39").Value = "  -4.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4626"
$ws.Range("E40").Value = "  -4.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8988"
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.50"
$ws.Range("E42").Value = "  -3.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.849"
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.365"
$ws.Range("E45").Value = "  -4.13%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.950"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1231"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.64"
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.8894"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4003"
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05724"
$ws.Range("E51").Value = "  -2.58%  "
